# Update ObjTables header metadata (new export timestamp + tableFormat="row")
# on sheet A1 (and A2 for the first sheet's file-level header) across all 26 tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("!!Compartment")
$ws.Unprotect()
$ws.Range("A1").Value = '!!!ObjTables objTablesVersion=''0.0.8'' date=''2020-03-09 15:33:14'''
$ws.Range("A2").Value = '!!ObjTables type=''Data'' id=''Compartment'' name=''Compartment'' date=''2020-03-09 15:33:14'' objTablesVersion=''0.0.8'' document=''TestModel layout'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Compound")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Compound'' name=''Compound'' date=''2020-03-09 15:33:14'' objTablesVersion=''0.0.8'' document=''TestModel layout'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Definition")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Definition'' name=''Definition'' date=''2020-03-09 15:33:14'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Enzyme")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Enzyme'' name=''Enzyme'' date=''2020-03-09 15:33:14'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!FbcObjective")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-03-09 15:33:14'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Gene")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Gene'' name=''Gene'' date=''2020-03-09 15:33:14'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Layout")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Layout'' name=''Layout'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' document=''TestModel layout'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Measurement")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Measurement'' name=''Measurement'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!PbConfig")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''PbConfig'' name=''PbConfig'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Position")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Position'' name=''Position'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Protein")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Protein'' name=''Protein'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Quantity")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Quantity'' name=''Quantity'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!QuantityInfo")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!QuantityMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Reaction")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Reaction'' name=''Reaction'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' document=''TestModel layout'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!ReactionStoichiometry")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Regulator")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Regulator'' name=''Regulator'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Relation")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Relation'' name=''Relation'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Relationship")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''Relationship'' name=''Relationship'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixColumn")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixOrdered")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixRow")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!StoichiometricMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!rxnconContingencyList")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item("!!rxnconReactionList")
$ws.Unprotect()
$ws.Range("A1").Value = '!!ObjTables type=''Data'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-03-09 15:33:15'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()
